$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 49, shifting rows 49-52 down to 50-53
$ws.Rows.Item(49).Insert()

# Populate the new row 49 with data (copy of original row 49 values with changes per diff)
$ws.Cells.Item(49, 1).Value = 9
$ws.Cells.Item(49, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(49, 3).Value = "Metropolitana"
$ws.Cells.Item(49, 4).Value = 45013
$ws.Cells.Item(49, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat
$ws.Cells.Item(49, 5).Value = 13
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100104
$ws.Cells.Item(49, 8).Value = "Frutos de pepita"
$ws.Cells.Item(49, 9).Value = 100104003
$ws.Cells.Item(49, 10).Value = "Membrillo"
$ws.Cells.Item(49, 11).Value = "Champion"
$ws.Cells.Item(49, 12).Value = "Primera"
$ws.Cells.Item(49, 13).Value = 450
$ws.Cells.Item(49, 14).Value = 11500
$ws.Cells.Item(49, 15).Value = 12000
$ws.Cells.Item(49, 16).Value = 11778
$ws.Cells.Item(49, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(49, 18).Value = "Provincia de Cachapoal"
$ws.Cells.Item(49, 19).Value = 785
$ws.Cells.Item(49, 20).Value = 15
